$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Phase 1: write every new formula first (no style tweaks yet) so
# that range-referencing formulas (the SUM totals) don't accidentally
# inherit a number format from an already-styled precedent cell.
# -----------------------------------------------------------------

# Group 1: rows 2-11 (10 items) -> markup price (J), totals (K11, L11)
$ws.Range("J2").Formula = '=ROUNDUP((E2*1.05), 2)'
$ws.Range("J3").Formula = '=ROUNDUP((E3-($J$2-$E$2)/9),2)'
$ws.Range("J4").Formula = '=ROUNDUP((E4-($J$2-$E$2)/9),2)'
$ws.Range("J5").Formula = '=ROUNDUP((E5-($J$2-$E$2)/9),2)'
$ws.Range("J6").Formula = '=ROUNDUP((E6-($J$2-$E$2)/9),2)'
$ws.Range("J7").Formula = '=ROUNDUP((E7-($J$2-$E$2)/9),2)'
$ws.Range("J8").Formula = '=ROUNDUP((E8-($J$2-$E$2)/9),2)'
$ws.Range("J9").Formula = '=ROUNDUP((E9-($J$2-$E$2)/9),2)'
$ws.Range("J10").Formula = '=ROUNDUP((E10-($J$2-$E$2)/9),2)'
$ws.Range("J11").Formula = '=ROUNDUP((E11-($J$2-$E$2)/9),2)'
$ws.Range("K11").Formula = '=SUM(J2:J11)'
$ws.Range("L11").Formula = '=SUM(E2:E11)'

# Group 2: rows 42-46 (5 items) -> markup price (J), totals (K46, L46)
$ws.Range("J42").Formula = '=ROUNDUP((E42-($J$46-$E$46)/9),2)'
$ws.Range("J43").Formula = '=ROUNDUP((E43-($J$46-$E$46)/9),2)'
$ws.Range("J44").Formula = '=ROUNDUP((E44-($J$46-$E$46)/9),2)'
$ws.Range("J45").Formula = '=ROUNDUP((E45-($J$46-$E$46)/9),2)'
$ws.Range("J46").Formula = '=ROUNDUP((E46*1.05), 2)'
$ws.Range("K46").Formula = '=SUM(J42:J46)'
$ws.Range("L46").Formula = '=SUM(E42:E46)'

# Group 3: rows 62-66 (5 items) -> markup price (J), totals (K66, L66)
$ws.Range("J62").Formula = '=ROUNDUP((E62*1.05), 2)'
$ws.Range("J63").Formula = '=ROUNDUP((E63-($J$62-$E$62)/9),2)'
$ws.Range("J64").Formula = '=ROUNDUP((E64-($J$62-$E$62)/9),2)'
$ws.Range("J65").Formula = '=ROUNDUP((E65-($J$62-$E$62)/9),2)'
$ws.Range("J66").Formula = '=ROUNDUP((E66-($J$62-$E$62)/9),2)'
$ws.Range("K66").Formula = '=SUM(J62:J66)'
$ws.Range("L66").Formula = '=SUM(E62:E66)'

# -----------------------------------------------------------------
# Phase 2: apply number formats / alignment now that every formula
# cell already exists, so nothing can "inherit" a format from a
# still-being-built precedent.
#   s=1 -> numFmt General, centered      (already existed)
#   s=2 -> numFmt 0.00,    centered      (already existed)
#   s=4 -> numFmt 0.00,    not centered  (new xf created by this edit)
# -----------------------------------------------------------------

# Group 1 -> all centered 0.00, including the K total; L total not centered
"J2","J3","J4","J5","J6","J7","J8","J9","J10","J11","K11" | ForEach-Object {
    $ws.Range($_).NumberFormat = "0.00"
    $ws.Range($_).HorizontalAlignment = -4108
}
$ws.Range("L11").NumberFormat = "0.00"

# Group 2 -> only J42:J45 and J46 are centered (J46 also needs 0.00); K46 stays default
"J42","J43","J44","J45" | ForEach-Object {
    $ws.Range($_).HorizontalAlignment = -4108
}
$ws.Range("J46").NumberFormat = "0.00"
$ws.Range("J46").HorizontalAlignment = -4108
$ws.Range("L46").NumberFormat = "0.00"

# Group 3 -> J62 is centered 0.00; J63:J66 centered only; K66 stays default
$ws.Range("J62").NumberFormat = "0.00"
$ws.Range("J62").HorizontalAlignment = -4108
"J63","J64","J65","J66" | ForEach-Object {
    $ws.Range($_).HorizontalAlignment = -4108
}
$ws.Range("L66").NumberFormat = "0.00"

# -----------------------------------------------------------------
# View state: selection moved to J13, scroll reset to top-left
# -----------------------------------------------------------------
$ws.Range("J13").Select()

$wb.Save()
